# DQ Report Query count.xlsx - weekly refresh
# Pulls forward last week's "latest" snapshot into the "previous" columns,
# records this week's new counts, and recomputes the diff column.
#
# date_latest_run   : 2025-08-19 (45888) -> 2025-08-26 (45895)
# date_previous_run : 2025-08-12 (45881) -> 2025-08-19 (45888)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLatestDate   = 45895
$newPreviousDate = 45888

# Row -> record_count_latest_run (C), record_count_previous_run (E), diff_since_last_run (F)
# $null means "leave that column blank" (used for rows with no previous-run history yet / anymore)
$rows = @(
    @{ Row = 2;  C = 3;   E = 4;   F = -1  },
    @{ Row = 3;  C = 53;  E = 45;  F = 8   },
    @{ Row = 4;  C = 4;   E = 6;   F = -2  },
    @{ Row = 5;  C = 23;  E = 35;  F = -12 },
    @{ Row = 6;  C = 38;  E = 30;  F = 8   },
    @{ Row = 7;  C = 42;  E = 44;  F = -2  },
    @{ Row = 8;  C = 6;   E = 6;   F = 0   },
    @{ Row = 9;  C = 2;   E = $null; F = $null },
    @{ Row = 10; C = 3;   E = 6;   F = -3  },
    @{ Row = 11; C = 6;   E = 6;   F = 0   },
    @{ Row = 12; C = 2;   E = 1;   F = 1   },
    @{ Row = 13; C = 7;   E = 5;   F = 2   },
    @{ Row = 14; C = 3;   E = 3;   F = 0   },
    @{ Row = 15; C = 16;  E = 16;  F = 0   },
    @{ Row = 16; C = 429; E = 493; F = -64 },
    @{ Row = 17; C = 19;  E = 21;  F = -2  },
    @{ Row = 18; C = 462; E = 445; F = 17  },
    @{ Row = 19; C = 79;  E = 65;  F = 14  },
    @{ Row = 20; C = 2;   E = 4;   F = -2  },
    @{ Row = 21; C = 141; E = 179; F = -38 },
    @{ Row = 22; C = 137; E = 173; F = -36 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("B$n").Value = $newLatestDate
    $ws.Range("C$n").Value = $r.C

    if ($null -eq $r.E) {
        $ws.Range("D$n").ClearContents()
        $ws.Range("E$n").ClearContents()
        $ws.Range("F$n").ClearContents()
    } else {
        $ws.Range("D$n").Value = $newPreviousDate
        $ws.Range("E$n").Value = $r.E
        $ws.Range("F$n").Value = $r.F
    }
}

# Row 23 ("06. COV Booster interval < 12 weeks") has no data for the new latest run yet;
# it only rolls its date_previous_run forward, record_count_previous_run stays as-is.
$ws.Range("D23").Value = $newPreviousDate
$ws.Range("E23").Value = 2

Write-Host "Updated DQ Report Query counts for latest run"
